$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bill_items")

# Insert a new column C ("item_description"), shifting item_type/quantity/price/total/tax
# one column to the right (D/E/F/G/H).
$ws.Columns.Item(3).Insert()

# New header cell
$ws.Range("C1").Value = "item_description"

# New data cell for the existing sample row
$ws.Range("C2").Value = "test bill item"

# The column Insert() operation copies the formatting of the adjacent column onto the
# newly inserted cells; the new description cell should have the default (no) style.
$ws.Range("C2").Style = "Normal"

# Match the new column's width (raw OOXML width of 16 characters).
$ws.Columns.Item(3).ColumnWidth = 15.1666666667

# Update the active selection to the new cell, as recorded in the saved view state.
$ws.Range("C2").Select() | Out-Null
